# Auto-generated edit script: updates market-price derived columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 393.33334
$ws.Range("I98").Value = 393.33334
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 393.33334
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 1104.66666
$ws.Range("N98").ClearContents()
$ws.Range("H100").Value = 65397.625
$ws.Range("I100").Value = 167738.33
$ws.Range("J100").Value = 3993.2
$ws.Range("K100").Value = 167738.33
$ws.Range("L100").Value = 3993.2
$ws.Range("M100").Value = -167197.33
$ws.Range("N100").Value = -5075.2
$ws.Range("H112").Value = 21444
$ws.Range("J112").Value = 26555
$ws.Range("L112").Value = 79665
$ws.Range("N112").Value = -81881
$ws.Range("H113").Value = 2184.4707
$ws.Range("I113").Value = 1541.5
$ws.Range("J113").Value = 3103
$ws.Range("K113").Value = 1541.5
$ws.Range("L113").Value = 3103
$ws.Range("M113").Value = 1712.5
$ws.Range("N113").Value = -9611
$ws.Range("H121").Value = 1335.2941
$ws.Range("J121").Value = 1485.7142
$ws.Range("L121").Value = 4457.142599999999
$ws.Range("N121").Value = -7951.142599999999
$ws.Range("H122").Value = 393.33334
$ws.Range("I122").Value = 393.33334
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1180.00002
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 1269.99998
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 42780
$ws.Range("J136").Value = 42780
$ws.Range("L136").Value = 42780
$ws.Range("N136").Value = -52980
$ws.Range("H138").Value = 4027.1836
$ws.Range("J138").Value = 8476.556
$ws.Range("L138").Value = 25429.668
$ws.Range("N138").Value = -35709.66800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 32260580
$ws.Range("I61").Value = 40002264
$ws.Range("J61").Value = 3571.3333
$ws.Range("K61").Value = 40002264
$ws.Range("L61").Value = 3571.3333
$ws.Range("M61").Value = -40002052
$ws.Range("N61").Value = -3995.3333
$ws.Range("H97").Value = 1159.52
$ws.Range("I97").Value = 886.2353000000001
$ws.Range("J97").Value = 1740.25
$ws.Range("K97").Value = 886.2353000000001
$ws.Range("L97").Value = 1740.25
$ws.Range("M97").Value = -390.2353000000001
$ws.Range("N97").Value = -2732.25
$ws.Range("H136").Value = 32260580
$ws.Range("I136").Value = 40002264
$ws.Range("J136").Value = 3571.3333
$ws.Range("K136").Value = 120006792
$ws.Range("L136").Value = 10713.9999
$ws.Range("M136").Value = -120004242
$ws.Range("N136").Value = -15813.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1312.7222
$ws.Range("I94").Value = 788.9
$ws.Range("J94").Value = 1967.5
$ws.Range("K94").Value = 788.9
$ws.Range("L94").Value = 1967.5
$ws.Range("M94").Value = -337.9
$ws.Range("N94").Value = -2869.5
$ws.Range("H107").Value = 2950
$ws.Range("I107").Value = 2900
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 2900
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -980
$ws.Range("N107").Value = -6840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 38.939392
$ws.Range("I7").Value = 31.277779
$ws.Range("J7").Value = 48.133335
$ws.Range("K7").Value = 31.277779
$ws.Range("L7").Value = 48.133335
$ws.Range("M7").Value = 81.722221
$ws.Range("N7").Value = -274.133335
$ws.Range("H107").Value = 506.17392
$ws.Range("I107").Value = 481.45
$ws.Range("J107").Value = 671
$ws.Range("K107").Value = 481.45
$ws.Range("L107").Value = 671
$ws.Range("M107").Value = 1438.55
$ws.Range("N107").Value = -4511

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H112").Value = 564.6667
$ws.Range("I112").Value = 564.6667
$ws.Range("K112").Value = 1694.0001
$ws.Range("M112").Value = -586.0001
$ws.Range("H140").Value = 57790.95
$ws.Range("I140").Value = 76268
$ws.Range("J140").Value = 2359.8
$ws.Range("K140").Value = 228804
$ws.Range("L140").Value = 7079.400000000001
$ws.Range("M140").Value = -223624
$ws.Range("N140").Value = -17439.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 583.25
$ws.Range("I97").Value = 545.36365
$ws.Range("K97").Value = 545.36365
$ws.Range("M97").Value = -49.36365000000001
$ws.Range("H107").Value = 251.72728
$ws.Range("I107").Value = 258.1
$ws.Range("J107").Value = 188
$ws.Range("K107").Value = 258.1
$ws.Range("L107").Value = 188
$ws.Range("M107").Value = 1661.9
$ws.Range("N107").Value = -4028
$ws.Range("H113").Value = 4029.1428
$ws.Range("I113").Value = 3967.9443
$ws.Range("K113").Value = 3967.9443
$ws.Range("M113").Value = -1797.9443
$ws.Range("H122").Value = 1581.6666
$ws.Range("I122").Value = 1400
$ws.Range("J122").Value = 1618
$ws.Range("K122").Value = 4200
$ws.Range("L122").Value = 4854
$ws.Range("M122").Value = -1750
$ws.Range("N122").Value = -9754

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1456.5333
$ws.Range("I7").Value = 1430.9
$ws.Range("J7").Value = 1507.8
$ws.Range("K7").Value = 1430.9
$ws.Range("L7").Value = 1507.8
$ws.Range("M7").Value = -1318.9
$ws.Range("N7").Value = -1731.8
$ws.Range("H69").Value = 48000
$ws.Range("J69").Value = 48000
$ws.Range("L69").Value = 48000
$ws.Range("N69").Value = -49622
$ws.Range("H72").Value = 48000
$ws.Range("J72").Value = 48000
$ws.Range("L72").Value = 144000
$ws.Range("N72").Value = -152112
$ws.Range("H93").Value = 1811.8334
$ws.Range("I93").Value = 1628.2858
$ws.Range("J93").Value = 2068.8
$ws.Range("K93").Value = 1628.2858
$ws.Range("L93").Value = 2068.8
$ws.Range("M93").Value = -380.2858000000001
$ws.Range("N93").Value = -4564.8
$ws.Range("H122").Value = 69469.53
$ws.Range("I122").Value = 168810.5
$ws.Range("J122").Value = 3242.2222
$ws.Range("K122").Value = 506431.5
$ws.Range("L122").Value = 9726.6666
$ws.Range("M122").Value = -503981.5
$ws.Range("N122").Value = -14626.6666
$ws.Range("H126").Value = 1456.5333
$ws.Range("I126").Value = 1430.9
$ws.Range("J126").Value = 1507.8
$ws.Range("K126").Value = 4292.700000000001
$ws.Range("L126").Value = 4523.4
$ws.Range("M126").Value = -1822.700000000001
$ws.Range("N126").Value = -9463.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2168.5833
$ws.Range("I96").Value = 1648.6
$ws.Range("J96").Value = 2540
$ws.Range("K96").Value = 1648.6
$ws.Range("L96").Value = 2540
$ws.Range("M96").Value = -275.5999999999999
$ws.Range("N96").Value = -5286
$ws.Range("H107").Value = 561.4286
$ws.Range("I107").Value = 532.5
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 1597.5
$ws.Range("L107").Value = 1800
$ws.Range("M107").Value = 322.5
$ws.Range("N107").Value = -5640
$ws.Range("H138").Value = 44429
$ws.Range("J138").Value = 44429
$ws.Range("L138").Value = 44429
$ws.Range("N138").Value = -54709
